# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-27 23:13:41
#
# Column G ("Recorded By") holds a comma-separated list of recorder names /
# emails (e.g. "System, someone@example.com"). This pass reverses the order
# of entries in that list for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colG = 7

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colG)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $text = [string]$current
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $reversed = $parts[($parts.Length - 1)..0]
    $newText = $reversed -join ", "

    $cell.Value = $newText
}
